# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" positioned right after "总计" (before
#    the existing "2022-Q3" sheet) with the fund-holding detail table.
# 2. Update the "总计" (totals) sheet to add a summary row for 2022-Q4,
#    shifting the existing 2022-Q3 / 2022-Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: build the new "2022-Q4" sheet.
#
# Duplicate the existing "2022-Q3" sheet (placing the copy right before
# it) so the new sheet inherits the exact same look & feel (bordered,
# bold, centred header row / index column), then drop all the sample
# rows but the first and overwrite it with the 2022-Q4 fund data.
# Final tab order: 总计, 2022-Q4, 2022-Q3, 2022-Q2
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet) | Out-Null
$q4Sheet = $wb.ActiveSheet
$q4Sheet.Name = "2022-Q4"

# Only one fund row is needed this quarter - drop the rest of the
# copied sample rows (rows 3-10 in the source sheet).
$q4Sheet.Range("A3:H10").Delete() | Out-Null

$q4Sheet.Range("A2").Value = 0

# Columns B, D, E, F, G hold text in the source data (e.g. "001900" with
# a leading zero, "0.13" / "62.72" as plain text) - force text format
# before assigning so the leading zero / decimal text isn't coerced to
# a number.
$q4Sheet.Range("B2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "001900"
$q4Sheet.Range("C2").Value = "诺安精选价值混合"
$q4Sheet.Range("D2").NumberFormat = "@"
$q4Sheet.Range("D2").Value = "0.13"
$q4Sheet.Range("E2").NumberFormat = "@"
$q4Sheet.Range("E2").Value = "62.72"
$q4Sheet.Range("F2").NumberFormat = "@"
$q4Sheet.Range("F2").Value = "1.75"
$q4Sheet.Range("G2").NumberFormat = "@"
$q4Sheet.Range("G2").Value = "0.0023"
$q4Sheet.Range("H2").Value = 8

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - insert a 2022-Q4 row above
# the existing 2022-Q3 / 2022-Q2 rows (which move down one row each).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Row 4 (new, previously unused) needs the same index-cell style
# ("s=2": bold + bordered + centred) as the existing A2/A3 cells -
# copy formatting from A3 before filling it in.
$totalSheet.Range("A3").Copy() | Out-Null
$totalSheet.Range("A4").PasteSpecial(-4122) | Out-Null    # xlPasteFormats

# Row 4: 2022-Q2 (shifted down from row 3)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.04

# Row 3: 2022-Q3 (shifted down from row 2)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 9
$totalSheet.Range("D3").Value = 0.08

# Row 2: 2022-Q4 (new)
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0

# Restore the originally-selected tab ("2022-Q2") now that all the
# sheet-shuffling/editing above is done (creating/copying sheets moves
# the active-sheet cursor around).
$wb.Worksheets.Item("2022-Q2").Activate()

Write-Output "2022-Q4 sheet added and totals updated"
